$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 28.666666
$ws.Range("I5").Value = 14.4
$ws.Range("K5").Value = 14.4
$ws.Range("M5").Value = 100.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 279.73
$ws.Range("I15").Value = 279.73
$ws.Range("K15").Value = 839.1900000000001
$ws.Range("M15").Value = -670.1900000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2889.48
$ws.Range("I138").Value = 2210.8
$ws.Range("J138").Value = 3907.5
$ws.Range("K138").Value = 6632.400000000001
$ws.Range("L138").Value = 11722.5
$ws.Range("M138").Value = -1492.400000000001
$ws.Range("N138").Value = -22002.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4500.6
$ws.Range("I32").Value = 4530.909
$ws.Range("J32").Value = 1500
$ws.Range("K32").Value = 4530.909
$ws.Range("L32").Value = 1500
$ws.Range("M32").Value = -4243.909
$ws.Range("N32").Value = -2074

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 232.72728
$ws.Range("J22").Value = 920
$ws.Range("L22").Value = 920
$ws.Range("N22").Value = -1266

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 56924.2
$ws.Range("I134").Value = 74952.266
$ws.Range("J134").Value = 2840
$ws.Range("K134").Value = 224856.798
$ws.Range("L134").Value = 8520
$ws.Range("M134").Value = -222321.798
$ws.Range("N134").Value = -13590

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 10615.2
$ws.Range("J50").Value = 10615.2
$ws.Range("L50").Value = 10615.2
$ws.Range("N50").Value = -11865.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 8760.875
$ws.Range("I51").Value = 1090
$ws.Range("J51").Value = 9856.714
$ws.Range("K51").Value = 1090
$ws.Range("L51").Value = 9856.714
$ws.Range("M51").Value = -354
$ws.Range("N51").Value = -11328.714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2544.4524
$ws.Range("I58").Value = 1098.6222
$ws.Range("J58").Value = 4212.718
$ws.Range("K58").Value = 1098.6222
$ws.Range("L58").Value = 4212.718
$ws.Range("M58").Value = -895.6222
$ws.Range("N58").Value = -4618.718

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 14458
$ws.Range("J59").Value = 15447.5
$ws.Range("L59").Value = 15447.5
$ws.Range("N59").Value = -17737.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 9368
$ws.Range("J60").Value = 10164
$ws.Range("L60").Value = 10164
$ws.Range("N60").Value = -11186

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 8760.875
$ws.Range("I61").Value = 1090
$ws.Range("J61").Value = 9856.714
$ws.Range("K61").Value = 1090
$ws.Range("L61").Value = 9856.714
$ws.Range("M61").Value = -742
$ws.Range("N61").Value = -10552.714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 17998.125
$ws.Range("J68").Value = 18531
$ws.Range("L68").Value = 18531
$ws.Range("N68").Value = -20029

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 17998.125
$ws.Range("J71").Value = 18531
$ws.Range("L71").Value = 55593
$ws.Range("N71").Value = -63081

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 14237.125
$ws.Range("J74").Value = 17418.666
$ws.Range("L74").Value = 17418.666
$ws.Range("N74").Value = -19166.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 14237.125
$ws.Range("J77").Value = 17418.666
$ws.Range("L77").Value = 52255.99800000001
$ws.Range("N77").Value = -60991.99800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 376028.5
$ws.Range("I99").Value = 168033.33
$ws.Range("K99").Value = 168033.33
$ws.Range("M99").Value = -166535.33

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 376028.5
$ws.Range("I126").Value = 168033.33
$ws.Range("K126").Value = 504099.99
$ws.Range("M126").Value = -501629.99

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2544.4524
$ws.Range("I136").Value = 1098.6222
$ws.Range("J136").Value = 4212.718
$ws.Range("K136").Value = 3295.8666
$ws.Range("L136").Value = 12638.154
$ws.Range("M136").Value = -745.8666000000003
$ws.Range("N136").Value = -17738.154

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 211.90909
$ws.Range("I98").Value = 233.16667
$ws.Range("J98").Value = 186.4
$ws.Range("K98").Value = 699.50001
$ws.Range("L98").Value = 559.2
$ws.Range("M98").Value = 798.49999
$ws.Range("N98").Value = -3555.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4213.4688
$ws.Range("J131").Value = 2567.9656
$ws.Range("L131").Value = 7703.8968
$ws.Range("N131").Value = -17783.8968

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4783.6665
$ws.Range("I70").Value = 4656.5713
$ws.Range("J70").Value = 4894.875
$ws.Range("K70").Value = 4656.5713
$ws.Range("L70").Value = 4894.875
$ws.Range("M70").Value = -4386.5713
$ws.Range("N70").Value = -5434.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4783.6665
$ws.Range("I73").Value = 4656.5713
$ws.Range("J73").Value = 4894.875
$ws.Range("K73").Value = 4656.5713
$ws.Range("L73").Value = 4894.875
$ws.Range("M73").Value = -3720.5713
$ws.Range("N73").Value = -6766.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2469.318
$ws.Range("I132").Value = 2070.182
$ws.Range("J132").Value = 2868.4546
$ws.Range("K132").Value = 6210.545999999999
$ws.Range("L132").Value = 8605.363799999999
$ws.Range("M132").Value = -3680.545999999999
$ws.Range("N132").Value = -13665.3638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 339.73334
$ws.Range("I22").Value = 377.33334
$ws.Range("J22").Value = 283.33334
$ws.Range("K22").Value = 377.33334
$ws.Range("L22").Value = 283.33334
$ws.Range("M22").Value = -82.33334000000002
$ws.Range("N22").Value = -873.33334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 339.73334
$ws.Range("I27").Value = 377.33334
$ws.Range("J27").Value = 283.33334
$ws.Range("K27").Value = 377.33334
$ws.Range("L27").Value = 283.33334
$ws.Range("M27").Value = -270.33334
$ws.Range("N27").Value = -497.33334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1533.2424
$ws.Range("I82").Value = 1336.826
$ws.Range("J82").Value = 1985
$ws.Range("K82").Value = 1336.826
$ws.Range("L82").Value = 1985
$ws.Range("M82").Value = -975.826
$ws.Range("N82").Value = -2707

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1533.2424
$ws.Range("I85").Value = 1336.826
$ws.Range("J85").Value = 1985
$ws.Range("K85").Value = 1336.826
$ws.Range("L85").Value = 1985
$ws.Range("M85").Value = -88.82600000000002
$ws.Range("N85").Value = -4481

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1135643.6
$ws.Range("I132").Value = 1595229.6
$ws.Range("J132").Value = 1998.2
$ws.Range("K132").Value = 4785688.800000001
$ws.Range("L132").Value = 5994.6
$ws.Range("M132").Value = -4783158.800000001
$ws.Range("N132").Value = -11054.6
